$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOB1004: Cálculo II (Requisito fraco)"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1004: Cálculo II (Requisito fraco)*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # The three paragraphs immediately following it (the empty paragraph,
    # the "Ver no Jupiter..." paragraph, and the "© 2020 ..." footer
    # paragraph) are removed, leaving the next (already-empty) paragraph
    # directly after "LOB1004: Cálculo II (Requisito fraco)".
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
